$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.121.47"
$ws.Range("E2").Value = "  +1.28%  "

Set-TextValue $ws.Range("D3") "3.932.16"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "493.60"
$ws.Range("E5").Value = "  +1.34%  "

Set-TextValue $ws.Range("D6") "147.46"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("E7").Value = "  -1.12%  "

$ws.Range("E8").Value = "  +0.03%  "

Set-TextValue $ws.Range("D9") "0.731"
$ws.Range("E9").Value = "  -0.38%  "

Set-TextValue $ws.Range("D10") "0.176"
$ws.Range("E10").Value = "  +4.28%  "

Set-TextValue $ws.Range("D11") "0.0000350"
$ws.Range("E11").Value = "  -0.72%  "

Set-TextValue $ws.Range("D12") "43.33"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("E13").Value = "  -2.18%  "

Set-TextValue $ws.Range("D14") "4.566.63"
$ws.Range("E14").Value = "  +0.37%  "

Set-TextValue $ws.Range("D15") "3.910.01"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("E16").Value = "  -3.56%  "

$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("E18").Value = "  +4.09%  "

Set-TextValue $ws.Range("D19") "19.86"
$ws.Range("E19").Value = "  -0.85%  "

Set-TextValue $ws.Range("D20") "69.211.65"
$ws.Range("E20").Value = "  +1.25%  "

Set-TextValue $ws.Range("D21") "439.13"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("E22").Value = "  +0.55%  "

Set-TextValue $ws.Range("D23") "14.53"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D24") "12.26"
$ws.Range("E24").Value = "  +9.55%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "88.63"
$ws.Range("E25").Value = "  +0.02%  "

Set-TextValue $ws.Range("D26") "3.81"
$ws.Range("E26").Value = "  +4.80%  "

Set-TextValue $ws.Range("D27") "11.14"
$ws.Range("E27").Value = "  -3.30%  "

Set-TextValue $ws.Range("D28") "37.10"
$ws.Range("E28").Value = "  -4.21%  "

Set-TextValue $ws.Range("D29") "5.66"
$ws.Range("E29").Value = "  -3.74%  "

Set-TextValue $ws.Range("D30") "703.72"
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("E31").Value = "  +0.20%  "

Set-TextValue $ws.Range("D32") "13.39"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("E33").Value = "  +0.07%  "

Set-TextValue $ws.Range("D34") "0.472"
$ws.Range("E34").Value = "  +18.86%  "

Set-TextValue $ws.Range("D35") "0.0₃0899"
$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D36") "6.12"
$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D37") "61.76"
$ws.Range("E37").Value = "  +3.80%  "

Set-TextValue $ws.Range("D38") "40.69"
$ws.Range("E38").Value = "  -2.97%  "

$ws.Range("E39").Value = "  +0.46%  "

Set-TextValue $ws.Range("D40") "0.996"
$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("E41").Value = "  +0.12%  "

Set-TextValue $ws.Range("D42") "0.0489"
$ws.Range("E42").Value = "  +1.38%  "

Set-TextValue $ws.Range("D43") "2.92"
$ws.Range("E43").Value = "  -2.10%  "

Set-TextValue $ws.Range("D44") "3.07"
$ws.Range("E44").Value = "  -3.86%  "

$ws.Range("E45").Value = "  +1.99%  "

$ws.Range("E46").Value = "  +0.44%  "

Set-TextValue $ws.Range("D47") "3.37"
$ws.Range("E47").Value = "  +7.18%  "

Set-TextValue $ws.Range("D48") "0.0₆0359"
$ws.Range("E48").Value = "  -0.58%  "

Set-TextValue $ws.Range("D49") "3.00"
$ws.Range("E49").Value = "  +5.99%  "

$ws.Range("E50").Value = "  -1.18%  "

Set-TextValue $ws.Range("D51") "144.34"
$ws.Range("E51").Value = "  -0.87%  "

